# Auto-generated Excel COM-interop script applying the scheduled-runner profit
# recompute described in the commit diff for Sheets/Marilith_Profits.xlsx.
# Each block updates the H:N (avg buy/sell price & profit) columns of one row
# in one of the 8 item-category worksheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2465.7144
$ws.Range("I2").Value = 712
$ws.Range("K2").Value = 712
$ws.Range("M2").Value = -599

$ws.Range("H61").Value = 1035
$ws.Range("J61").Value = 1522.5
$ws.Range("L61").Value = 4567.5
$ws.Range("N61").Value = -4911.5

$ws.Range("H74").Value = 168984.67
$ws.Range("I74").Value = 3500
$ws.Range("K74").Value = 3500
$ws.Range("M74").Value = -2564

$ws.Range("H77").Value = 168984.67
$ws.Range("I77").Value = 3500
$ws.Range("K77").Value = 17500
$ws.Range("M77").Value = -12820

$ws.Range("H132").Value = 1846.25
$ws.Range("I132").Value = 1586.1904
$ws.Range("K132").Value = 4758.5712
$ws.Range("M132").Value = -2228.5712

$ws.Range("H138").Value = 2219.3333
$ws.Range("I138").Value = 1136.4
$ws.Range("J138").Value = 2992.8572
$ws.Range("K138").Value = 3409.2
$ws.Range("L138").Value = 8978.571599999999
$ws.Range("M138").Value = 1730.8
$ws.Range("N138").Value = -19258.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 31513
$ws.Range("I33").Value = 30026
$ws.Range("K33").Value = 30026
$ws.Range("M33").Value = -29697

$ws.Range("H40").Value = 30515.5
$ws.Range("J40").Value = 30515.5
$ws.Range("L40").Value = 30515.5
$ws.Range("N40").Value = -30867.5

$ws.Range("H45").Value = 1424.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1424.5
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1424.5
$ws.Range("N45").Value = -2178.5
$ws.Range("M45").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 799.6667
$ws.Range("I74").Value = 699.5
$ws.Range("K74").Value = 699.5
$ws.Range("M74").Value = 174.5

$ws.Range("H77").Value = 799.6667
$ws.Range("I77").Value = 699.5
$ws.Range("K77").Value = 3497.5
$ws.Range("M77").Value = 870.5

$ws.Range("H132").Value = 1699.7142
$ws.Range("I132").Value = 1699.7142
$ws.Range("K132").Value = 5099.142599999999
$ws.Range("M132").Value = -2569.142599999999

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 149.25
$ws.Range("I22").Value = 131.66667
$ws.Range("J22").Value = 202
$ws.Range("K22").Value = 131.66667
$ws.Range("L22").Value = 202
$ws.Range("M22").Value = 41.33332999999999
$ws.Range("N22").Value = -548

$ws.Range("H29").Value = 2339
$ws.Range("J29").Value = 5018
$ws.Range("L29").Value = 5018
$ws.Range("N29").Value = -5596

$ws.Range("H39").Value = 6000
$ws.Range("J39").Value = 6000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6778

$ws.Range("H134").Value = 4320.8823
$ws.Range("I134").Value = 3778.5
$ws.Range("K134").Value = 11335.5
$ws.Range("M134").Value = -8800.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 39.2
$ws.Range("I7").Value = 33.142857
$ws.Range("J7").Value = 53.333332
$ws.Range("K7").Value = 33.142857
$ws.Range("L7").Value = 53.333332
$ws.Range("M7").Value = 79.85714300000001
$ws.Range("N7").Value = -279.333332

$ws.Range("H12").Value = 5848.75
$ws.Range("J12").Value = 5848.75
$ws.Range("L12").Value = 5848.75
$ws.Range("N12").Value = -6188.75

$ws.Range("H58").Value = 3224.5
$ws.Range("I58").Value = 2299.3333
$ws.Range("K58").Value = 2299.3333
$ws.Range("M58").Value = -2096.3333

$ws.Range("H105").Value = 1637.5
$ws.Range("I105").Value = 1637.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1637.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 109.5
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 616.55554
$ws.Range("I107").Value = 507.14285
$ws.Range("K107").Value = 507.14285
$ws.Range("M107").Value = 1412.85715

$ws.Range("H136").Value = 3224.5
$ws.Range("I136").Value = 2299.3333
$ws.Range("K136").Value = 6897.999899999999
$ws.Range("M136").Value = -4347.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99710
$ws.Range("J37").Value = 99710
$ws.Range("L37").Value = 299130
$ws.Range("N37").Value = -299354

$ws.Range("H40").Value = 457.16666
$ws.Range("I40").Value = 148.8
$ws.Range("J40").Value = 1999
$ws.Range("K40").Value = 595.2
$ws.Range("L40").Value = 7996
$ws.Range("M40").Value = -526.2
$ws.Range("N40").Value = -8134

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H97").Value = 275.5
$ws.Range("I97").Value = 149.5
$ws.Range("J97").Value = 401.5
$ws.Range("K97").Value = 448.5
$ws.Range("L97").Value = 1204.5
$ws.Range("M97").Value = 47.5
$ws.Range("N97").Value = -2196.5

$ws.Range("H113").Value = 787.2
$ws.Range("I113").Value = 546.5
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 1639.5
$ws.Range("L113").Value = 5250
$ws.Range("M113").Value = 530.5
$ws.Range("N113").Value = -9590

$ws.Range("H136").Value = 3100
$ws.Range("I136").Value = 3100
$ws.Range("K136").Value = 9300
$ws.Range("M136").Value = -4200

$ws.Range("H139").Value = 671.875
$ws.Range("I139").Value = 671.875
$ws.Range("K139").Value = 2015.625
$ws.Range("M139").Value = 3124.375

$ws.Range("H140").Value = 1367.5
$ws.Range("I140").Value = 1367.5
$ws.Range("K140").Value = 4102.5
$ws.Range("M140").Value = 1077.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 222224060
$ws.Range("I70").Value = 222224060
$ws.Range("K70").Value = 222224060
$ws.Range("M70").Value = -222223790

$ws.Range("H73").Value = 222224060
$ws.Range("I73").Value = 222224060
$ws.Range("K73").Value = 222224060
$ws.Range("M73").Value = -222223124

$ws.Range("H80").Value = 1769.8
$ws.Range("J80").Value = 1999.8334
$ws.Range("L80").Value = 1999.8334
$ws.Range("N80").Value = -3995.8334

$ws.Range("H83").Value = 1769.8
$ws.Range("J83").Value = 1999.8334
$ws.Range("L83").Value = 9999.166999999999
$ws.Range("N83").Value = -19983.167

$ws.Range("H122").Value = 1432.1111
$ws.Range("I122").Value = 1414.2858
$ws.Range("J122").Value = 1494.5
$ws.Range("K122").Value = 4242.857400000001
$ws.Range("L122").Value = 4483.5
$ws.Range("M122").Value = -1792.857400000001
$ws.Range("N122").Value = -9383.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 873.5
$ws.Range("I16").Value = 873.5
$ws.Range("K16").Value = 873.5
$ws.Range("M16").Value = -703.5

$ws.Range("H23").Value = 3000
$ws.Range("J23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("N23").Value = -3460

$ws.Range("H46").Value = 3750
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812

$ws.Range("H56").Value = 14999.5
$ws.Range("I56").Value = 14999
$ws.Range("K56").Value = 14999
$ws.Range("M56").Value = -14308

$ws.Range("H81").Value = 39493.5
$ws.Range("J81").Value = 39493.5
$ws.Range("L81").Value = 39493.5
$ws.Range("N81").Value = -41489.5

$ws.Range("H84").Value = 39493.5
$ws.Range("J84").Value = 39493.5
$ws.Range("L84").Value = 118480.5
$ws.Range("N84").Value = -128464.5

$ws.Range("H134").Value = 95000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 95000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 95000
$ws.Range("N134").Value = -105140
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 30000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 30000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 30000
$ws.Range("N61").Value = -30584
$ws.Range("M61").ClearContents()

$ws.Range("H119").Value = 15000
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 15000
$ws.Range("N119").Value = -24676
